$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("B241").Value = "좌표공간의 주어진 점을 `$x y`$ 평면에 대하여 대칭이동한 점의 좌표를 구합니다."
$ws.Range("B242").Value = "좌표공간의 주어진 점을 `$y z`$ 평면에 대하여 대칭이동한 점의 좌표를 구합니다."
$ws.Range("B243").Value = "좌표공간의 두 점 사이의 거리를 구합니다."

$ws.Range("A242").Value = "z0002"
$ws.Range("A243").Value = "z0003"
$ws.Range("A244").Value = "z0004"
$ws.Range("A245").Value = "z0005"
$ws.Range("A246").Value = "z0006"
$ws.Range("A247").Value = "z0007"
$ws.Range("A248").Value = "z0008"
$ws.Range("A249").Value = "z0009"

$ws.Range("B244").Value = "주어진 초점의 좌표와 일치하도록 쌍곡선식의 미정계수를 구합니다."
$ws.Range("B245").Value = "쌍곡선의 방정식에서 주축의 길이를 계산합니다."
$ws.Range("B246").Value = "좌표평면 위의 두 직선의 방향벡터를 구합니다."
$ws.Range("B247").Value = "두 벡터의 내적을 이용해서 `$cos\theta`$를 구합니다."
$ws.Range("B248").Value = "타원 위의 점과 두 초점의 거리의 합을 구합니다."
$ws.Range("B249").Value = "사각형을 높이가 같은 두 개의 삼각형으로 나눠서 생각합니다."

$ws.Range("A250").Value = "z0010"
$ws.Range("B250").Value = "두 삼각형의 넓이의 합이 사각형의 넓이라는 사실로 방정식을 세웁니다."
$ws.Range("A251").Value = "z0011"
$ws.Range("B251").Value = "원의 반지름의 길이를 구합니다."

$ws.Range("B251").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 227
$win.ScrollColumn = 1
